# Auto-generated edit script: refresh market-price derived columns (H-N)
# per scheduled-runner update, matching the authoritative diff.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 375471.44
$ws.Range("I15").Value = 375471.44
$ws.Range("K15").Value = 1126414.32
$ws.Range("M15").Value = -1126245.32
$ws.Range("H33").Value = 747.1818
$ws.Range("I33").Value = 696.6667
$ws.Range("K33").Value = 696.6667
$ws.Range("M33").Value = -467.6667
$ws.Range("H74").Value = 12919.3
$ws.Range("I74").Value = 6922
$ws.Range("K74").Value = 6922
$ws.Range("M74").Value = -5986
$ws.Range("H77").Value = 12919.3
$ws.Range("I77").Value = 6922
$ws.Range("K77").Value = 34610
$ws.Range("M77").Value = -29930
$ws.Range("H132").Value = 3563.6428
$ws.Range("I132").Value = 2847.28
$ws.Range("J132").Value = 9533.333000000001
$ws.Range("K132").Value = 8541.84
$ws.Range("L132").Value = 28599.999
$ws.Range("M132").Value = -6011.84
$ws.Range("N132").Value = -33659.999

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4112.5405
$ws.Range("J2").Value = 6068.2856
$ws.Range("L2").Value = 6068.2856
$ws.Range("N2").Value = -6294.2856
$ws.Range("H45").Value = 1792.1562
$ws.Range("I45").Value = 1787.4
$ws.Range("J45").Value = 1809.1428
$ws.Range("K45").Value = 1787.4
$ws.Range("L45").Value = 1809.1428
$ws.Range("M45").Value = -1410.4
$ws.Range("N45").Value = -2563.1428
$ws.Range("H63").Value = 3824.6316
$ws.Range("I63").Value = 2438.5557
$ws.Range("J63").Value = 5072.1
$ws.Range("K63").Value = 2438.5557
$ws.Range("L63").Value = 5072.1
$ws.Range("M63").Value = -1752.5557
$ws.Range("N63").Value = -6444.1
$ws.Range("H66").Value = 3824.6316
$ws.Range("I66").Value = 2438.5557
$ws.Range("J66").Value = 5072.1
$ws.Range("K66").Value = 12192.7785
$ws.Range("L66").Value = 25360.5
$ws.Range("M66").Value = -8760.7785
$ws.Range("N66").Value = -32224.5
$ws.Range("H116").Value = 4112.5405
$ws.Range("J116").Value = 6068.2856
$ws.Range("L116").Value = 6068.2856
$ws.Range("N116").Value = -10656.2856
$ws.Range("H122").Value = 3664.2632
$ws.Range("I122").Value = 3511.8
$ws.Range("J122").Value = 3833.6667
$ws.Range("K122").Value = 10535.4
$ws.Range("L122").Value = 11501.0001
$ws.Range("M122").Value = -8085.400000000001
$ws.Range("N122").Value = -16401.0001

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4112.5405
$ws.Range("J3").Value = 6068.2856
$ws.Range("L3").Value = 6068.2856
$ws.Range("N3").Value = -6296.2856
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H86").Value = 29949
$ws.Range("I86").Value = 31248.75
$ws.Range("J86").Value = 24750
$ws.Range("K86").Value = 31248.75
$ws.Range("L86").Value = 24750
$ws.Range("M86").Value = -30125.75
$ws.Range("N86").Value = -26996
$ws.Range("H89").Value = 29949
$ws.Range("I89").Value = 31248.75
$ws.Range("J89").Value = 24750
$ws.Range("K89").Value = 156243.75
$ws.Range("L89").Value = 123750
$ws.Range("M89").Value = -150627.75
$ws.Range("N89").Value = -134982
$ws.Range("H94").Value = 700.63635
$ws.Range("I94").Value = 789.52
$ws.Range("J94").Value = 422.875
$ws.Range("K94").Value = 789.52
$ws.Range("L94").Value = 422.875
$ws.Range("M94").Value = -338.52
$ws.Range("N94").Value = -1324.875

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8467.571
$ws.Range("I16").Value = 7637.5
$ws.Range("J16").Value = 8799.6
$ws.Range("K16").Value = 7637.5
$ws.Range("L16").Value = 8799.6
$ws.Range("M16").Value = -7350.5
$ws.Range("N16").Value = -9373.6
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H31").Value = 3309.4666
$ws.Range("I31").Value = 1388.8
$ws.Range("K31").Value = 1388.8
$ws.Range("M31").Value = -1093.8
$ws.Range("H34").Value = 3309.4666
$ws.Range("I34").Value = 1388.8
$ws.Range("K34").Value = 1388.8
$ws.Range("M34").Value = -1186.8
$ws.Range("H58").Value = 33338504
$ws.Range("I58").Value = 45457596
$ws.Range("K58").Value = 45457596
$ws.Range("M58").Value = -45457393
$ws.Range("H62").Value = 12541.5
$ws.Range("I62").Value = 11516.667
$ws.Range("J62").Value = 13566.333
$ws.Range("K62").Value = 11516.667
$ws.Range("L62").Value = 13566.333
$ws.Range("M62").Value = -10892.667
$ws.Range("N62").Value = -14814.333
$ws.Range("H65").Value = 12541.5
$ws.Range("I65").Value = 11516.667
$ws.Range("J65").Value = 13566.333
$ws.Range("K65").Value = 57583.335
$ws.Range("L65").Value = 67831.66500000001
$ws.Range("M65").Value = -54463.335
$ws.Range("N65").Value = -74071.66500000001
$ws.Range("H99").Value = 8551174
$ws.Range("I99").Value = 11114927
$ws.Range("K99").Value = 11114927
$ws.Range("M99").Value = -11113429
$ws.Range("H105").Value = 8024.2856
$ws.Range("J105").Value = 9747.166999999999
$ws.Range("L105").Value = 9747.166999999999
$ws.Range("N105").Value = -13241.167
$ws.Range("H107").Value = 669.38464
$ws.Range("I107").Value = 588.5263
$ws.Range("J107").Value = 888.8570999999999
$ws.Range("K107").Value = 588.5263
$ws.Range("L107").Value = 888.8570999999999
$ws.Range("M107").Value = 1331.4737
$ws.Range("N107").Value = -4728.8571
$ws.Range("H113").Value = 8467.571
$ws.Range("I113").Value = 7637.5
$ws.Range("J113").Value = 8799.6
$ws.Range("K113").Value = 7637.5
$ws.Range("L113").Value = 8799.6
$ws.Range("M113").Value = -5467.5
$ws.Range("N113").Value = -13139.6
$ws.Range("H122").Value = 4499.375
$ws.Range("I122").Value = 4318
$ws.Range("J122").Value = 4608.2
$ws.Range("K122").Value = 12954
$ws.Range("L122").Value = 13824.6
$ws.Range("M122").Value = -10504
$ws.Range("N122").Value = -18724.6
$ws.Range("H126").Value = 8551174
$ws.Range("I126").Value = 11114927
$ws.Range("K126").Value = 33344781
$ws.Range("M126").Value = -33342311
$ws.Range("H136").Value = 33338504
$ws.Range("I136").Value = 45457596
$ws.Range("K136").Value = 136372788
$ws.Range("M136").Value = -136370238

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 207640.83
$ws.Range("I4").Value = 319.5
$ws.Range("K4").Value = 958.5
$ws.Range("M4").Value = -846.5
$ws.Range("H98").Value = 1496.2
$ws.Range("J98").Value = 1120.125
$ws.Range("L98").Value = 3360.375
$ws.Range("N98").Value = -6356.375

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1500
$ws.Range("J13").Value = 1500
$ws.Range("L13").Value = 1500
$ws.Range("N13").Value = -1778
$ws.Range("H23").Value = 3711
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 3711
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 3711
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -4157
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H80").Value = 18894.889
$ws.Range("J80").Value = 24107
$ws.Range("L80").Value = 24107
$ws.Range("N80").Value = -26103
$ws.Range("H83").Value = 18894.889
$ws.Range("J83").Value = 24107
$ws.Range("L83").Value = 120535
$ws.Range("N83").Value = -130519
$ws.Range("H102").Value = 4781.6523
$ws.Range("I102").Value = 4411.6113
$ws.Range("K102").Value = 4411.6113
$ws.Range("M102").Value = -2789.6113

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4333.25
$ws.Range("I40").Value = 3990
$ws.Range("J40").Value = 4600.222
$ws.Range("K40").Value = 3990
$ws.Range("L40").Value = 4600.222
$ws.Range("M40").Value = -3854
$ws.Range("N40").Value = -4872.222
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H136").Value = 3572.15
$ws.Range("I136").Value = 1447.5
$ws.Range("J136").Value = 4482.7144
$ws.Range("K136").Value = 4342.5
$ws.Range("L136").Value = 13448.1432
$ws.Range("M136").Value = -1792.5
$ws.Range("N136").Value = -18548.1432

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 1749.25
$ws.Range("J31").Value = 1749.25
$ws.Range("L31").Value = 1749.25
$ws.Range("N31").Value = -2445.25
$ws.Range("H126").Value = 3511.12
$ws.Range("I126").Value = 1999.1428
$ws.Range("J126").Value = 5435.4546
$ws.Range("K126").Value = 5997.428400000001
$ws.Range("L126").Value = 16306.3638
$ws.Range("M126").Value = -3527.428400000001
$ws.Range("N126").Value = -21246.3638
$ws.Range("H136").Value = 18561020
$ws.Range("I136").Value = 21784208
$ws.Range("K136").Value = 65352624
$ws.Range("M136").Value = -65350074

Write-Host "Applied Odin_Profits market-data refresh."